$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: semantic identifiers change from iaest-measure: to iaest-dimension:
# for the three dimension columns (Edad, Sexo, Mes y año)
$ws.Range("A3").Value = "iaest-dimension:edad-grupos-quinquenales"
$ws.Range("F3").Value = "iaest-dimension:sexo"
$ws.Range("G3").Value = "iaest-dimension:mes-y-ano"

# Row 4: those same dimension columns switch from "medida" to "dim"
$ws.Range("A4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("G4").Value = "dim"

# Row 5: those same dimension columns switch from "xsd:string" to "skos:Concept"
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("F5").Value = "skos:Concept"

# Row 6 (new): mapping file references for the two dimensions that now
# resolve against an external concept mapping workbook.
# Write the values first, then copy the formatting from row 5 so the new
# cells pick up the same cell style used throughout the sheet.
$ws.Range("A6").Value = "mapping-edad-grupos-quinquenales.xlsx"
$ws.Range("F6").Value = "mapping-sexo.xlsx"

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122)
